$wb = $excel.ActiveWorkbook

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 6811.9414
$ws_ARM.Range("I2").Value = 885
$ws_ARM.Range("J2").Value = 10044.818
$ws_ARM.Range("K2").Value = 885
$ws_ARM.Range("L2").Value = 10044.818
$ws_ARM.Range("M2").Value = -772
$ws_ARM.Range("N2").Value = -10270.818

$ws_ARM.Range("H32").Value = 3910.8872
$ws_ARM.Range("I32").Value = 2657.6924
$ws_ARM.Range("J32").Value = 17487.166
$ws_ARM.Range("K32").Value = 2657.6924
$ws_ARM.Range("L32").Value = 17487.166
$ws_ARM.Range("M32").Value = -2370.6924
$ws_ARM.Range("N32").Value = -18061.166

$ws_ARM.Range("H45").Value = 2244.2104
$ws_ARM.Range("I45").Value = 2001.7333
$ws_ARM.Range("J45").Value = 3153.5
$ws_ARM.Range("K45").Value = 2001.7333
$ws_ARM.Range("L45").Value = 3153.5
$ws_ARM.Range("M45").Value = -1624.7333
$ws_ARM.Range("N45").Value = -3907.5

$ws_ARM.Range("H61").Value = 2433.25
$ws_ARM.Range("I61").Value = 929.9
$ws_ARM.Range("J61").Value = 9950
$ws_ARM.Range("K61").Value = 929.9
$ws_ARM.Range("L61").Value = 9950
$ws_ARM.Range("M61").Value = -717.9
$ws_ARM.Range("N61").Value = -10374

$ws_ARM.Range("H74").Value = 1042.8679
$ws_ARM.Range("I74").Value = 1030.4688
$ws_ARM.Range("J74").Value = 1061.762
$ws_ARM.Range("K74").Value = 1030.4688
$ws_ARM.Range("L74").Value = 1061.762
$ws_ARM.Range("M74").Value = -156.4688000000001
$ws_ARM.Range("N74").Value = -2809.762

$ws_ARM.Range("H77").Value = 1042.8679
$ws_ARM.Range("I77").Value = 1030.4688
$ws_ARM.Range("J77").Value = 1061.762
$ws_ARM.Range("K77").Value = 5152.344000000001
$ws_ARM.Range("L77").Value = 5308.809999999999
$ws_ARM.Range("M77").Value = -784.344000000001
$ws_ARM.Range("N77").Value = -14044.81

$ws_ARM.Range("H116").Value = 6811.9414
$ws_ARM.Range("I116").Value = 885
$ws_ARM.Range("J116").Value = 10044.818
$ws_ARM.Range("K116").Value = 885
$ws_ARM.Range("L116").Value = 10044.818
$ws_ARM.Range("M116").Value = 1409
$ws_ARM.Range("N116").Value = -14632.818

$ws_ARM.Range("H122").Value = 1222.0714
$ws_ARM.Range("I122").Value = 1070.9
$ws_ARM.Range("K122").Value = 3212.7
$ws_ARM.Range("M122").Value = -762.7000000000003

$ws_ARM.Range("H132").Value = 3652.5293
$ws_ARM.Range("I132").Value = 1339.3334
$ws_ARM.Range("J132").Value = 7893.3887
$ws_ARM.Range("K132").Value = 4018.0002
$ws_ARM.Range("L132").Value = 23680.1661
$ws_ARM.Range("M132").Value = -1488.0002
$ws_ARM.Range("N132").Value = -28740.1661

$ws_ARM.Range("H136").Value = 2433.25
$ws_ARM.Range("I136").Value = 929.9
$ws_ARM.Range("J136").Value = 9950
$ws_ARM.Range("K136").Value = 2789.7
$ws_ARM.Range("L136").Value = 29850
$ws_ARM.Range("M136").Value = -239.6999999999998
$ws_ARM.Range("N136").Value = -34950

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 6811.9414
$ws_BSM.Range("I3").Value = 885
$ws_BSM.Range("J3").Value = 10044.818
$ws_BSM.Range("K3").Value = 885
$ws_BSM.Range("L3").Value = 10044.818
$ws_BSM.Range("M3").Value = -771
$ws_BSM.Range("N3").Value = -10272.818

$ws_BSM.Range("H134").Value = 3123
$ws_BSM.Range("I134").Value = 1869.3572
$ws_BSM.Range("J134").Value = 4585.5835
$ws_BSM.Range("K134").Value = 5608.071599999999
$ws_BSM.Range("L134").Value = 13756.7505
$ws_BSM.Range("M134").Value = -3073.071599999999
$ws_BSM.Range("N134").Value = -18826.7505

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 15050
$ws_CRP.Range("I31").Value = 0
$ws_CRP.Range("J31").Value = 15050
$ws_CRP.Range("K31").Value = 0
$ws_CRP.Range("L31").Value = 15050
$ws_CRP.Range("M31").Value = ""
$ws_CRP.Range("N31").Value = -15640

$ws_CRP.Range("H34").Value = 15050
$ws_CRP.Range("I34").Value = 0
$ws_CRP.Range("J34").Value = 15050
$ws_CRP.Range("K34").Value = 0
$ws_CRP.Range("L34").Value = 15050
$ws_CRP.Range("M34").Value = ""
$ws_CRP.Range("N34").Value = -15454

$ws_CRP.Range("H107").Value = 513.7778
$ws_CRP.Range("I107").Value = 230.14285
$ws_CRP.Range("J107").Value = 1506.5
$ws_CRP.Range("K107").Value = 230.14285
$ws_CRP.Range("L107").Value = 1506.5
$ws_CRP.Range("M107").Value = 1689.85715
$ws_CRP.Range("N107").Value = -5346.5

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 3697.8386
$ws_CUL.Range("I5").Value = 378.18182
$ws_CUL.Range("J5").Value = 11812.556
$ws_CUL.Range("K5").Value = 1134.54546
$ws_CUL.Range("L5").Value = 35437.66800000001
$ws_CUL.Range("M5").Value = -1022.54546
$ws_CUL.Range("N5").Value = -35661.66800000001

$ws_CUL.Range("H33").Value = 1698.6364
$ws_CUL.Range("I33").Value = 429.9375
$ws_CUL.Range("J33").Value = 5081.8335
$ws_CUL.Range("K33").Value = 2579.625
$ws_CUL.Range("L33").Value = 30491.001
$ws_CUL.Range("M33").Value = -2296.625
$ws_CUL.Range("N33").Value = -31057.001

$ws_CUL.Range("H44").Value = 83334420
$ws_CUL.Range("I44").Value = 142857380
$ws_CUL.Range("J44").Value = 45456164
$ws_CUL.Range("K44").Value = 428572140
$ws_CUL.Range("L44").Value = 136368492
$ws_CUL.Range("M44").Value = -428571742
$ws_CUL.Range("N44").Value = -136369288

$ws_CUL.Range("H69").Value = 25001478
$ws_CUL.Range("J69").Value = 28573068
$ws_CUL.Range("L69").Value = 85719204
$ws_CUL.Range("N69").Value = -85720826

$ws_CUL.Range("H72").Value = 25001478
$ws_CUL.Range("J72").Value = 28573068
$ws_CUL.Range("L72").Value = 257157612
$ws_CUL.Range("N72").Value = -257165724

$ws_CUL.Range("H80").Value = 2633.3333
$ws_CUL.Range("J80").Value = 2623.5293
$ws_CUL.Range("L80").Value = 7870.5879
$ws_CUL.Range("N80").Value = -9742.5879

$ws_CUL.Range("H83").Value = 2633.3333
$ws_CUL.Range("J83").Value = 2623.5293
$ws_CUL.Range("L83").Value = 23611.7637
$ws_CUL.Range("N83").Value = -32971.7637

$ws_CUL.Range("H107").Value = 163.57895
$ws_CUL.Range("I107").Value = 63
$ws_CUL.Range("J107").Value = 175.41176
$ws_CUL.Range("K107").Value = 189
$ws_CUL.Range("L107").Value = 526.23528
$ws_CUL.Range("M107").Value = 1731
$ws_CUL.Range("N107").Value = -4366.23528

$ws_CUL.Range("H132").Value = 842.875
$ws_CUL.Range("I132").Value = 697.5
$ws_CUL.Range("J132").Value = 988.25
$ws_CUL.Range("K132").Value = 6277.5
$ws_CUL.Range("L132").Value = 8894.25
$ws_CUL.Range("M132").Value = -3747.5
$ws_CUL.Range("N132").Value = -13954.25

$ws_CUL.Range("H135").Value = 3697.8386
$ws_CUL.Range("I135").Value = 378.18182
$ws_CUL.Range("J135").Value = 11812.556
$ws_CUL.Range("K135").Value = 3403.63638
$ws_CUL.Range("L135").Value = 106313.004
$ws_CUL.Range("M135").Value = -868.6363799999999
$ws_CUL.Range("N135").Value = -111383.004
